$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): update 想去人数 (F column) for the two events
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 1300
$wsExhibit.Range("F3").Value = 2817

# Sheet "全部类型" (sheet4): same two events appear here, update accordingly
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 1300
$wsAll.Range("F4").Value = 2817
